# "Update Stückliste (dual opamp), pdf vom Schema, Feedback PH"
#
# The old single op-amp "LMV796MF/NOPB" (Digikey "LMV796MF/NOPBCT-ND") used
# in position 26 ("Verstärker", row 28) is replaced by the dual op-amp
# "LMV797MM/NOPB" (Digikey "LMV797MM/NOPBCT-ND") at the new unit price of
# 1.82 CHF (was 1.46 CHF). The "Kosten" column is a formula (H*I) so it -
# and the grand total in J34 - recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bauteil (part number) and Bestell Nr. (order number) for position 26
$ws.Range("D28").Value = "LMV797MM/NOPB"
$ws.Range("F28").Value = "LMV797MM/NOPBCT-ND"

# Updated unit price ("Stückpreis") for the new part
$ws.Range("I28").Value = 1.82

# Leave the cursor where the editor last left it
$ws.Range("I29").Select()
